$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 4
$ws.Range("H4").Value = 1699.75
$ws.Range("I4").Value = 1766.3334
$ws.Range("K4").Value = 1766.3334
$ws.Range("M4").Value = -1652.3334

# Row 58
$ws.Range("H58").Value = 1194
$ws.Range("I58").Value = 1000
$ws.Range("J58").Value = 1388
$ws.Range("K58").Value = 3000
$ws.Range("L58").Value = 4164
$ws.Range("M58").Value = -2850
$ws.Range("N58").Value = -4464

# Row 61
$ws.Range("H61").Value = 400
$ws.Range("I61").Value = 400
$ws.Range("K61").Value = 1200
$ws.Range("M61").Value = -1028

# Row 95
$ws.Range("H95").Value = 50000
$ws.Range("J95").Value = 50000
$ws.Range("L95").Value = 50000
$ws.Range("N95").Value = -55492

# Row 137
$ws.Range("H137").Value = 33595.676
$ws.Range("I137").Value = 1038.4706
$ws.Range("K137").Value = 3115.4118
$ws.Range("M137").Value = -565.4118000000003

# Row 138
$ws.Range("H138").Value = 1749.6063
$ws.Range("I138").Value = 1421.2616
$ws.Range("J138").Value = 2485.5518
$ws.Range("K138").Value = 4263.7848
$ws.Range("L138").Value = 7456.655400000001
$ws.Range("M138").Value = 876.2151999999996
$ws.Range("N138").Value = -17736.6554

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 5113.42
$ws.Range("I32").Value = 4670.884
$ws.Range("J32").Value = 7831.857
$ws.Range("K32").Value = 4670.884
$ws.Range("L32").Value = 7831.857
$ws.Range("M32").Value = -4383.884
$ws.Range("N32").Value = -8405.857

# Row 61
$ws.Range("H61").Value = 2076.6
$ws.Range("I61").Value = 1404.3462
$ws.Range("K61").Value = 1404.3462
$ws.Range("M61").Value = -1192.3462

# Row 74
$ws.Range("H74").Value = 2581.5833
$ws.Range("I74").Value = 2340.2856
$ws.Range("K74").Value = 2340.2856
$ws.Range("M74").Value = -1466.2856

# Row 77
$ws.Range("H77").Value = 2581.5833
$ws.Range("I77").Value = 2340.2856
$ws.Range("K77").Value = 11701.428
$ws.Range("M77").Value = -7333.428

# Row 136
$ws.Range("H136").Value = 2076.6
$ws.Range("I136").Value = 1404.3462
$ws.Range("K136").Value = 4213.0386
$ws.Range("M136").Value = -1663.0386

$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 3699.6592
$ws.Range("I134").Value = 3730.6052
$ws.Range("J134").Value = 3503.6667
$ws.Range("K134").Value = 11191.8156
$ws.Range("L134").Value = 10511.0001
$ws.Range("M134").Value = -8656.8156
$ws.Range("N134").Value = -15581.0001

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2442.3
$ws.Range("I31").Value = 1869.3334
$ws.Range("K31").Value = 1869.3334
$ws.Range("M31").Value = -1574.3334

# Row 34
$ws.Range("H34").Value = 2442.3
$ws.Range("I34").Value = 1869.3334
$ws.Range("K34").Value = 1869.3334
$ws.Range("M34").Value = -1667.3334

# Row 74
$ws.Range("H74").Value = 34999.668
$ws.Range("J74").Value = 34999.668
$ws.Range("L74").Value = 34999.668
$ws.Range("N74").Value = -36747.668

# Row 77
$ws.Range("H77").Value = 34999.668
$ws.Range("J77").Value = 34999.668
$ws.Range("L77").Value = 104999.004
$ws.Range("N77").Value = -113735.004

# Row 132
$ws.Range("H132").Value = 1680.151
$ws.Range("I132").Value = 1207.8975
$ws.Range("K132").Value = 3623.6925
$ws.Range("M132").Value = -1093.6925

# Row 134
$ws.Range("H134").Value = 1919.4186
$ws.Range("I134").Value = 1775.5143
$ws.Range("K134").Value = 5326.5429
$ws.Range("M134").Value = -2791.5429

$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 288350.7
$ws.Range("I4").Value = 123.75
$ws.Range("J4").Value = 749513.8
$ws.Range("K4").Value = 371.25
$ws.Range("L4").Value = 2248541.4
$ws.Range("M4").Value = -259.25
$ws.Range("N4").Value = -2248765.4

# Row 12
$ws.Range("H12").Value = 132.17647
$ws.Range("J12").Value = 171.54546
$ws.Range("L12").Value = 514.6363799999999
$ws.Range("N12").Value = -860.6363799999999

# Row 81
$ws.Range("H81").Value = 1274.5
$ws.Range("I81").Value = 549
$ws.Range("J81").Value = 2000
$ws.Range("K81").Value = 1647
$ws.Range("L81").Value = 6000
$ws.Range("M81").Value = -524
$ws.Range("N81").Value = -8246

# Row 84
$ws.Range("H84").Value = 1274.5
$ws.Range("I84").Value = 549
$ws.Range("J84").Value = 2000
$ws.Range("K84").Value = 4941
$ws.Range("L84").Value = 18000
$ws.Range("M84").Value = 675
$ws.Range("N84").Value = -29232

# Row 87
$ws.Range("H87").Value = 11935.333
$ws.Range("I87").Value = 5403
$ws.Range("J87").Value = 25000
$ws.Range("K87").Value = 16209
$ws.Range("L87").Value = 75000
$ws.Range("M87").Value = -14961
$ws.Range("N87").Value = -77496

# Row 90
$ws.Range("H90").Value = 11935.333
$ws.Range("I90").Value = 5403
$ws.Range("J90").Value = 25000
$ws.Range("K90").Value = 48627
$ws.Range("L90").Value = 225000
$ws.Range("M90").Value = -42387
$ws.Range("N90").Value = -237480

$ws = $wb.Worksheets.Item("GSM")
# Row 46
$ws.Range("H46").Value = 22975
$ws.Range("J46").Value = 22975
$ws.Range("L46").Value = 22975
$ws.Range("N46").Value = -23287

# Row 132
$ws.Range("H132").Value = 803425.9
$ws.Range("I132").Value = 1283824.1
$ws.Range("K132").Value = 3851472.3
$ws.Range("M132").Value = -3848942.3

# Row 134
$ws.Range("H134").Value = 27795.834
$ws.Range("J134").Value = 27795.834
$ws.Range("L134").Value = 83387.50199999999
$ws.Range("N134").Value = -88457.50199999999

$ws = $wb.Worksheets.Item("LTW")
# Row 54
$ws.Range("H54").Value = 39977
$ws.Range("J54").Value = 39977
$ws.Range("L54").Value = 39977
$ws.Range("N54").Value = -41265

# Row 55
$ws.Range("H55").Value = 536.8421
$ws.Range("I55").Value = 509
$ws.Range("J55").Value = 584.5714
$ws.Range("K55").Value = 509
$ws.Range("L55").Value = 584.5714
$ws.Range("M55").Value = -336
$ws.Range("N55").Value = -930.5714

# Row 64
$ws.Range("H64").Value = 502999.5
$ws.Range("I64").Value = 999999
$ws.Range("J64").Value = 6000
$ws.Range("K64").Value = 999999
$ws.Range("L64").Value = 6000
$ws.Range("M64").Value = -999774
$ws.Range("N64").Value = -6450

# Row 67
$ws.Range("H67").Value = 502999.5
$ws.Range("I67").Value = 999999
$ws.Range("J67").Value = 6000
$ws.Range("K67").Value = 999999
$ws.Range("L67").Value = 6000
$ws.Range("M67").Value = -999219
$ws.Range("N67").Value = -7560

# Row 68
$ws.Range("H68").Value = 1489.75
$ws.Range("J68").Value = 1798.9
$ws.Range("L68").Value = 1798.9
$ws.Range("N68").Value = -3296.9

# Row 71
$ws.Range("H71").Value = 1489.75
$ws.Range("J71").Value = 1798.9
$ws.Range("L71").Value = 8994.5
$ws.Range("N71").Value = -16482.5

# Row 132
$ws.Range("H132").Value = 1949.326
$ws.Range("I132").Value = 1577.2667
$ws.Range("K132").Value = 4731.800099999999
$ws.Range("M132").Value = -2201.800099999999

# Row 51
$ws.Range("H51").Value = 29977
$ws.Range("J51").Value = 29977
$ws.Range("L51").Value = 29977
$ws.Range("N51").Value = -30997

$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 5285.909
$ws.Range("I62").Value = 4855.75
$ws.Range("K62").Value = 4855.75
$ws.Range("M62").Value = -4231.75

# Row 65
$ws.Range("H65").Value = 5285.909
$ws.Range("I65").Value = 4855.75
$ws.Range("K65").Value = 24278.75
$ws.Range("M65").Value = -21158.75

# Row 80
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()

# Row 82
$ws.Range("H82").Value = 10000
$ws.Range("J82").Value = 10000
$ws.Range("L82").Value = 10000
$ws.Range("N82").Value = -10766

# Row 83
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()

# Row 85
$ws.Range("H85").Value = 10000
$ws.Range("J85").Value = 10000
$ws.Range("L85").Value = 10000
$ws.Range("N85").Value = -12652

# Row 95
$ws.Range("H95").Value = 99998
$ws.Range("J95").Value = 99998
$ws.Range("L95").Value = 99998
$ws.Range("N95").Value = -105490

# Row 122
$ws.Range("H122").Value = 41202.152
$ws.Range("I122").Value = 54927.414
$ws.Range("K122").Value = 164782.242
$ws.Range("M122").Value = -162332.242

# Row 126
$ws.Range("H126").Value = 4547.028
$ws.Range("I126").Value = 4640.6665
$ws.Range("J126").Value = 4266.1113
$ws.Range("K126").Value = 13921.9995
$ws.Range("L126").Value = 12798.3339
$ws.Range("M126").Value = -11451.9995
$ws.Range("N126").Value = -17738.3339

# Row 136
$ws.Range("H136").Value = 10289973
$ws.Range("I136").Value = 12628156
$ws.Range("J136").Value = 1967.5
$ws.Range("K136").Value = 37884468
$ws.Range("L136").Value = 5902.5
$ws.Range("M136").Value = -37881918
$ws.Range("N136").Value = -11002.5
